# Gantt chart workbook update:
#  - A new tracking row (date 2024-06-11 / 45478, Open=20, Closed=47) is
#    inserted just above the closing "thick border" rows of the table on
#    'Main Board', pushing the table's bottom border rows (and the chart
#    anchored below them) down by one row.
#  - The area chart's two series are extended to include the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Board")

# ---------------------------------------------------------------------
# 1. Insert a new row 34 (pushes the old row 34/35 border rows to 35/36,
#    and leaves rows 1-33 - and the shared D/E formulas living in them -
#    untouched).
# ---------------------------------------------------------------------
$ws.Range("A34:E34").EntireRow.Insert()

# Pick up the formatting (number formats / borders / styles) of the row
# directly above (the last data row) so the new row matches the existing
# table styling exactly.
$ws.Range("A33:E33").Copy()
$ws.Range("A34:E34").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Fill in the new row's data + formulas.
# ---------------------------------------------------------------------
$ws.Cells.Item(34, 1).Value = 45478
$ws.Cells.Item(34, 2).Value = 20
$ws.Cells.Item(34, 3).Value = 47
$ws.Cells.Item(34, 4).Formula = "=B34+C34"
$ws.Cells.Item(34, 5).Formula = "=C34/D34"

# Match the authored selection (cell D34, the newly-added "Days" value).
$ws.Range("D34").Select()

# ---------------------------------------------------------------------
# 3. Extend the area chart's two series ('Closed' / 'Open') so they pick
#    up the new row of data.
# ---------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

$closedSeries = $chart.SeriesCollection().Item(1)
$closedSeries.Formula = "=SERIES('Main Board'!`$C`$1,'Main Board'!`$A`$2:`$A`$34,'Main Board'!`$C`$2:`$C`$34,1)"

$openSeries = $chart.SeriesCollection().Item(2)
$openSeries.Formula = "=SERIES('Main Board'!`$B`$1,'Main Board'!`$A`$2:`$A`$34,'Main Board'!`$B`$2:`$B`$34,2)"

# ---------------------------------------------------------------------
# 4. The chart object itself is anchored below the table; since a row was
#    inserted above it, move it down by exactly the height of that new
#    row so it keeps sitting right under the (now one-row-lower) table.
# ---------------------------------------------------------------------
$fromRowOffPts = 48985 / 12700
$toRowOffPts = 16328 / 12700

$topTotal = 0
for ($r = 1; $r -le 36; $r++) {
    $topTotal += $ws.Rows.Item($r).Height
}
$bottomTotal = 0
for ($r = 1; $r -le 51; $r++) {
    $bottomTotal += $ws.Rows.Item($r).Height
}

$newTop = $topTotal + $fromRowOffPts
$newBottom = $bottomTotal + $toRowOffPts

$co.Top = $newTop
$co.Height = $newBottom - $newTop
